# Apply the changes described by the commit "prepared for simulator B testing 4"
# to the overview_generation sheet (new rows 6-9) and update the selections
# on the overview_generation and input_variables sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("overview_generation")
$wsInput    = $wb.Worksheets.Item("input_variables")

# ---------------------------------------------------------------------------
# Add the four new data rows (6-9) to overview_generation, re-using the
# formatting of row 5 (wrap-text cells A:E, date cell C, hyperlink cell F).
# ---------------------------------------------------------------------------

$wsOverview.Range("A5:F5").Copy() | Out-Null
$wsOverview.Range("A6:F9").PasteSpecial(-4122) | Out-Null
$wsOverview.Application.CutCopyMode = $false

# Hyperlinks for column F, rows 6-9 (added before the text values are set so
# that our own cell text - not the hyperlink's default display text - ends
# up in the cell).
$wsOverview.Hyperlinks.Add($wsOverview.Range("F6"), "", "initial_parameter_values!A1", "", "initial_parameters") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("F7"), "", "initial_parameter_values!A1", "", "initial_parameters") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("F8"), "", "initial_parameter_values!A1", "", "initial_parameters") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("F9"), "", "initial_parameter_values!A1", "", "initial_parameters") | Out-Null

# Restore the formatting the hyperlink creation might have altered, and set
# row heights to match the wrapped two-line text rows.
$wsOverview.Range("A5:F5").Copy() | Out-Null
$wsOverview.Range("A6:F9").PasteSpecial(-4122) | Out-Null
$wsOverview.Application.CutCopyMode = $false
$wsOverview.Rows.Item(6).RowHeight = 43.5
$wsOverview.Rows.Item(7).RowHeight = 43.5
$wsOverview.Rows.Item(8).RowHeight = 43.5
$wsOverview.Rows.Item(9).RowHeight = 43.5

# Text values are entered in the same order the original author used, which
# controls the order new entries are appended to the shared-string table.
$wsOverview.Range("A6").Value = "20240507_simulator_B_generated"
$wsOverview.Range("A7").Value = "20240520_simulator_A_generated"
$wsOverview.Range("A8").Value = "20240521_simulator_B_generated"
$wsOverview.Range("A9").Value = "20240522_simulator_A_generated"

$wsOverview.Range("D7").Value = "4th intervall for simulator A"
$wsOverview.Range("D9").Value = "5th intervall for simulator A"
$wsOverview.Range("D6").Value = "3rd time"
$wsOverview.Range("D8").Value = "4th time"

$wsOverview.Range("E6").Value = "tested different parameter combinations with = data generation. Kept constant simulation interval."
$wsOverview.Range("E7").Value = "tested different parameter combinations with = data generation. Kept constant simulation interval."
$wsOverview.Range("E8").Value = "tested different parameter combinations with = data generation. Kept constant simulation interval."
$wsOverview.Range("E9").Value = "tested different parameter combinations with = data generation. Kept constant simulation interval."

$wsOverview.Range("B6").Value = "B"
$wsOverview.Range("B7").Value = "A"
$wsOverview.Range("B8").Value = "B"
$wsOverview.Range("B9").Value = "A"

$wsOverview.Range("C6").Value = 45419
$wsOverview.Range("C7").Value = 45432
$wsOverview.Range("C8").Value = 45433
$wsOverview.Range("C9").Value = 45434

$wsOverview.Range("F6").Value = "initial_parameter_values"
$wsOverview.Range("F7").Value = "initial_parameter_values"
$wsOverview.Range("F8").Value = "initial_parameter_values"
$wsOverview.Range("F9").Value = "initial_parameter_values"

# ---------------------------------------------------------------------------
# Update view/selection state.
# ---------------------------------------------------------------------------

# overview_generation: selection moves to D5.
$wsOverview.Activate() | Out-Null
$wsOverview.Range("D5").Select() | Out-Null

# input_variables: stays the active/selected tab, selection moves to C13.
$wsInput.Activate() | Out-Null
$wsInput.Range("C13").Select() | Out-Null
